$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51; this shifts existing rows 51-135 down to 52-136,
# automatically extending the used range / dimension to R136 and carrying the
# date-format style from the row above into the new row's D cell.
$ws.Rows.Item(51).Insert()

# Populate the new row 51 with the new data record.
$ws.Cells.Item(51, 1).Value = 7
$ws.Cells.Item(51, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(51, 3).Value = "Ñuble"
$ws.Cells.Item(51, 4).Value = 44495
$ws.Cells.Item(51, 5).Value = 16
$ws.Cells.Item(51, 6).Value = 100112017
$ws.Cells.Item(51, 7).Value = "Apio"
$ws.Cells.Item(51, 8).Value = "Americana (o)"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 100
$ws.Cells.Item(51, 11).Value = 8000
$ws.Cells.Item(51, 12).Value = 8500
$ws.Cells.Item(51, 13).Value = 8250
$ws.Cells.Item(51, 14).Value = "$/docena de matas"
$ws.Cells.Item(51, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(51, 16).Value = 1375
$ws.Cells.Item(51, 17).Value = 6
$ws.Cells.Item(51, 18).Value = "Hortaliza"
